$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new ticker row right after the last used row (A83 -> A84)
$ws.Range("A84").Value = "GRT-USD"
